$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.508.32"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -4.39%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.261.20"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.82%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "554.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.29%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.590"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.81%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.253.96"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.185"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.97%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.585"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.62%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.09"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -6.70%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000264"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "632.59"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.55"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.55%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.788.98"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.64%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.503.23"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.13%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.79"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.57%  "
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.116"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.267.65"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.91%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.901"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.92%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.63"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "106.12"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +8.72%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -7.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.96"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.79%  "
$ws.Range("E27").Value = "  -5.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.49"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.35%  "
$ws.Range("E29").Value = "  -5.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.33"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.95"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.50%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.26"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.95%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "10.99"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.83%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "543.93"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +12.63%  "
$ws.Range("E35").Value = "  -2.55%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "57.03"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.82%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.645.72"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.46%  "
$ws.Range("E39").Value = "  -1.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.130"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0₃0718"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.58%  "
$ws.Range("E42").Value = "  -4.43%  "
$ws.Range("E43").Value = "  -3.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "32.23"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.335"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.27"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.92%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0414"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.61"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.60%  "
$ws.Range("E49").Value = "  -3.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("E51").Value = "  +1.51%  "
